$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of row 2 and row 3 for the columns that differ between
# the two species records (A, B, E, F, G, H, Q, R). Also row 3's empty
# "Bestamningsmetod" (AF) cell moves to row 2.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $tmp = $ws.Range($addr2).Value2
    $ws.Range($addr2).Value2 = $ws.Range($addr3).Value2
    $ws.Range($addr3).Value2 = $tmp
}

# I ("Antal") holds a numeric-looking label ("30") that must stay text, so
# force the destination cell to Text format before writing it - otherwise
# Excel's normal type inference would store it as the number 30.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value2 = "30"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value2 = ""

# AF2 becomes an (empty) cell, same as AF3 was before the edit, while AF3
# loses its (empty) cell entry. Touch the number format so the now-empty
# cells stay materialised instead of being dropped as fully blank.
$ws.Range("AF2").NumberFormat = "@"
$ws.Range("AF2").Value2 = ""
$ws.Range("AF3").ClearContents()
